# Split the single run of text in the Title, Author and Abstract
# paragraphs into one run per word (with the separating single spaces
# also becoming their own runs), leaving the text content identical.

$d = $word.ActiveDocument

function Split-IntoWordRuns([string]$text) {
    # Split on single spaces, keeping the spaces themselves as separate
    # tokens (there are no double-spaces in the target text, so a plain
    # split on " " followed by re-inserting a lone-space token between
    # every pair of words reproduces the original string exactly).
    $words = $text -split ' '
    $tokens = New-Object System.Collections.Generic.List[string]
    for ($i = 0; $i -lt $words.Count; $i++) {
        if ($i -gt 0) { [void]$tokens.Add(' ') }
        [void]$tokens.Add($words[$i])
    }
    return $tokens
}

function Set-ParagraphRunsPerWord($paragraph) {
    $range = $paragraph.Range
    $fullText = $paragraph.Range.Text
    # Paragraph.Range.Text includes the trailing paragraph mark; strip it
    # (and any cell/row marks) before tokenising, we only want the
    # visible text.
    $fullText = $fullText -replace "[\r\a\f]+$", ""

    $styleName = $paragraph.Style.NameLocal

    $tokens = Split-IntoWordRuns $fullText

    $runsXml = New-Object System.Text.StringBuilder
    foreach ($tok in $tokens) {
        $escaped = $tok.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
        [void]$runsXml.Append('<w:r><w:t xml:space="preserve">')
        [void]$runsXml.Append($escaped)
        [void]$runsXml.Append('</w:t></w:r>')
    }

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p><w:pPr><w:pStyle w:val="' + $styleName + '"/></w:pPr>' +
        $runsXml.ToString() +
        '</w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    [void]$range.InsertXML($xml)
}

# Title: "Questions: The scalar product"
Set-ParagraphRunsPerWord $d.Paragraphs.Item(1)

# Author: "Ritwik Anand"
Set-ParagraphRunsPerWord $d.Paragraphs.Item(2)

# Abstract: "A selection of questions for the study guide on the scalar product"
Set-ParagraphRunsPerWord $d.Paragraphs.Item(4)
